# Update "想去人数" (want-to-go count) figures that were refreshed by the
# site's data export (gh-pages output regenerated at commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 49
$ws.Range("F7").Value = 1111
$ws.Range("F11").Value = 7898
$ws.Range("F13").Value = 9262
$ws.Range("F30").Value = 14
$ws.Range("F35").Value = 324
$ws.Range("F37").Value = 891
$ws.Range("F41").Value = 400

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F20").Value = 343

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 49
$ws.Range("F7").Value = 1111
$ws.Range("F13").Value = 7898
$ws.Range("F14").Value = 9262
$ws.Range("F22").Value = 14
$ws.Range("F28").Value = 324
$ws.Range("F31").Value = 891
$ws.Range("F36").Value = 400
$ws.Range("F47").Value = 343
